# Split the "Source PubMed ID (PMID) or doi" column into two separate
# columns - "Source PubMed ID (PMID)" (numeric) and "Source DOI" - on both
# the "Score Development Samples" and "Evaluation Sample Sets" sheets.

$wb = $excel.ActiveWorkbook

# --- Score Development Samples -------------------------------------------
$ws = $wb.Worksheets.Item("Score Development Samples")

# Insert a new column at P (16); it pushes the old P/Q ("Cohort(s)" /
# "Additional Sample/Cohort Information") one slot to the right and
# inherits the header style from the surrounding header row.
$ws.Columns.Item(16).Insert()

# Rename the original (now PMID-only) header and give the new column its
# own "Source DOI" header.
$ws.Range("O1").Value = "Source PubMed ID (PMID)"
$ws.Range("P1").Value = "Source DOI"

# The PMID values that used to live in O as text are now plain numbers.
$ws.Range("O2").Value = 10000012
$ws.Range("O3").Value = 10000013

# --- Evaluation Sample Sets ------------------------------------------------
$ws2 = $wb.Worksheets.Item("Evaluation Sample Sets")

$ws2.Columns.Item(16).Insert()

$ws2.Range("O1").Value = "Source PubMed ID (PMID)"
$ws2.Range("P1").Value = "Source DOI"

# This sheet's sample row had no PMID, but it did carry a DOI value - move
# it into the new "Source DOI" column.
$ws2.Range("P2").Value = "10.2021/pgs.1003"
